# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
# timestamps on the per-language sheets to reflect a freshly-generated report.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-23 21:18:45"
$zhcn.Range("H2").Value = "2016-03-23 21:19:10"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-23 21:18:49"
$dede.Range("H2").Value = "2016-03-23 21:19:17"
